$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 58 (shifts existing rows 58..153 down to 59..154)
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new record
$ws.Cells.Item(58, 1).Value = 5
$ws.Cells.Item(58, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(58, 3).Value = "Maule"
$ws.Cells.Item(58, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(58, 5).Value = 7
$ws.Cells.Item(58, 6).Value = 100112021
$ws.Cells.Item(58, 7).Value = "Ají"
$ws.Cells.Item(58, 8).Value = "Americana (o)"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 150
$ws.Cells.Item(58, 11).Value = 15000
$ws.Cells.Item(58, 12).Value = 15000
$ws.Cells.Item(58, 13).Value = 15000
$ws.Cells.Item(58, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(58, 15).Value = "Región del Maule"
$ws.Cells.Item(58, 16).Value = 1071
$ws.Cells.Item(58, 17).Value = 14
$ws.Cells.Item(58, 18).Value = "Hortaliza"
